# Update the two existing filename entries and append a new "Flu09" row,
# matching the data/upload state described in the commit "set data and upload".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap/relabel the existing filenames in column B.
$ws.Range("B2").Value = "cytokine_data.xlsx"
$ws.Range("B3").Value = "patient_data.xlsx"

# Add the new row (A4/B4) that extends the table.
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Flu09"

# Give A4 the same number/border/font style as the other index cells (A2:A3)
# by copying formats only from A3, so no extra cell-style is introduced.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null
